$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.820.47"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "'2.219.59"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'292.20"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").Value = "'86.22"
$ws.Range("E6").Value = "  +5.19%  "
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("D10").Value = "'30.80"
$ws.Range("E10").Value = "  +6.27%  "
$ws.Range("E11").Value = "  +1.98%  "
$ws.Range("D12").Value = "'47.21"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").Value = "'2.562.74"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "'2.218.22"
$ws.Range("D18").Value = "'0.732"
$ws.Range("E18").Value = "  +3.00%  "
$ws.Range("D19").Value = "'39.788.19"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").Value = "'0.0₃0882"
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("D21").Value = "'11.07"
$ws.Range("E21").Value = "  +7.76%  "
$ws.Range("D22").Value = "'5.80"
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("D23").Value = "'65.65"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "'236.42"
$ws.Range("E24").Value = "  +4.85%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +2.77%  "
$ws.Range("D27").Value = "'1.84"
$ws.Range("E27").Value = "  +3.03%  "
$ws.Range("D28").Value = "'22.76"
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("D29").Value = "'2.20"
$ws.Range("E29").Value = "  +4.47%  "
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").Value = "'32.89"
$ws.Range("E31").Value = "  +4.20%  "
$ws.Range("D32").Value = "'151.63"
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").Value = "'4.94"
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("E35").Value = "  +3.98%  "
$ws.Range("D36").Value = "'2.37"
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("E37").Value = "  +6.97%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'15.86"
$ws.Range("E39").Value = "  +4.73%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.0994"
$ws.Range("E40").Value = "  +3.04%  "
$ws.Range("E41").Value = "  +3.93%  "
$ws.Range("E42").Value = "  +4.92%  "
$ws.Range("D43").Value = "'2.066.63"
$ws.Range("E43").Value = "  +9.20%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'2.10"
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("E45").Value = "  +3.98%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'17.82"
$ws.Range("E46").Value = "  +11.85%  "
$ws.Range("D47").Value = "'9.93"
$ws.Range("E47").Value = "  +11.14%  "
$ws.Range("D48").Value = "'2.60"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "'2.435.91"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "'71.37"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "'89.00"
$ws.Range("E51").Value = "  +2.70%  "
